$d = $word.ActiveDocument
$total = 0

# --- Paragraph 1 (PT Programa) ---
$r = $d.Content.Find.Execute("ão e condutividade hidráulica)Movimentos gravitacionais de m", $true, $false, $false, $false, $false, $true, 1, $false, "ão e condutividade hidráulica)^lMovimentos gravitacionais de m", 2)
if (-not $r) { Write-Output "FAILED: ão e condutividade hidráulica)Movimentos gravitacionais de m" }
$total += 1
$r = $d.Content.Find.Execute("s e técnicas de estabilização Processos erosivos: conceitos,", $true, $false, $false, $false, $false, $true, 1, $false, "s e técnicas de estabilização ^lProcessos erosivos: conceitos,", 2)
if (-not $r) { Write-Output "FAILED: s e técnicas de estabilização Processos erosivos: conceitos," }
$total += 1
$r = $d.Content.Find.Execute(" processos erosivos antrópicosColapso e subsidência", $true, $false, $false, $false, $false, $true, 1, $false, " processos erosivos antrópicos^lColapso e subsidência", 2)
if (-not $r) { Write-Output "FAILED:  processos erosivos antrópicosColapso e subsidência" }
$total += 1
$r = $d.Content.Find.Execute("Colapso e subsidênciaAterros de resíduos (seleção d", $true, $false, $false, $false, $false, $true, 1, $false, "Colapso e subsidência^lAterros de resíduos (seleção d", 2)
if (-not $r) { Write-Output "FAILED: Colapso e subsidênciaAterros de resíduos (seleção d" }
$total += 1
$r = $d.Content.Find.Execute("tivos, ABNT 8419 e ABNT 15849)Levantamento de estudos e proj", $true, $false, $false, $false, $false, $true, 1, $false, "tivos, ABNT 8419 e ABNT 15849)^lLevantamento de estudos e proj", 2)
if (-not $r) { Write-Output "FAILED: tivos, ABNT 8419 e ABNT 15849)Levantamento de estudos e proj" }
$total += 1
$r = $d.Content.Find.Execute("aplicados a Estudos AmbientaisAplicação de geossintéticos em", $true, $false, $false, $false, $false, $true, 1, $false, "aplicados a Estudos Ambientais^lAplicação de geossintéticos em", 2)
if (-not $r) { Write-Output "FAILED: aplicados a Estudos AmbientaisAplicação de geossintéticos em" }
$total += 1
$r = $d.Content.Find.Execute("sintéticos em obras ambientaisAulas práticas: parâmetros de ", $true, $false, $false, $false, $false, $true, 1, $false, "sintéticos em obras ambientais^lAulas práticas: parâmetros de ", 2)
if (-not $r) { Write-Output "FAILED: sintéticos em obras ambientaisAulas práticas: parâmetros de " }
$total += 1
$r = $d.Content.Find.Execute("âmetros de mecânica dos solos A disciplina pode contar com v", $true, $false, $false, $false, $false, $true, 1, $false, "âmetros de mecânica dos solos ^lA disciplina pode contar com v", 2)
if (-not $r) { Write-Output "FAILED: âmetros de mecânica dos solos A disciplina pode contar com v" }
$total += 1

# --- Paragraph 2 (EN Programa) ---
$r = $d.Content.Find.Execute("on and hydraulic conductivity)Gravitational mass movements: ", $true, $false, $false, $false, $false, $true, 1, $false, "on and hydraulic conductivity)^lGravitational mass movements: ", 2)
if (-not $r) { Write-Output "FAILED: on and hydraulic conductivity)Gravitational mass movements: " }
$total += 1
$r = $d.Content.Find.Execute("s and stabilization techniquesErosive processes: concepts, t", $true, $false, $false, $false, $false, $true, 1, $false, "s and stabilization techniques^lErosive processes: concepts, t", 2)
if (-not $r) { Write-Output "FAILED: s and stabilization techniquesErosive processes: concepts, t" }
$total += 1
$r = $d.Content.Find.Execute("nthropogenic erosion processesCollapse and subsidence", $true, $false, $false, $false, $false, $true, 1, $false, "nthropogenic erosion processes^lCollapse and subsidence", 2)
if (-not $r) { Write-Output "FAILED: nthropogenic erosion processesCollapse and subsidence" }
$total += 1
$r = $d.Content.Find.Execute("Collapse and subsidenceWaste landfills (selection of ", $true, $false, $false, $false, $false, $true, 1, $false, "Collapse and subsidence^lWaste landfills (selection of ", 2)
if (-not $r) { Write-Output "FAILED: Collapse and subsidenceWaste landfills (selection of " }
$total += 1
$r = $d.Content.Find.Execute("ods, ABNT 8419 and ABNT 15849)Survey of geotechnical studies", $true, $false, $false, $false, $false, $true, 1, $false, "ods, ABNT 8419 and ABNT 15849)^lSurvey of geotechnical studies", 2)
if (-not $r) { Write-Output "FAILED: ods, ABNT 8419 and ABNT 15849)Survey of geotechnical studies" }
$total += 1
$r = $d.Content.Find.Execute("plied to Environmental StudiesApplication of geosynthetics i", $true, $false, $false, $false, $false, $true, 1, $false, "plied to Environmental Studies^lApplication of geosynthetics i", 2)
if (-not $r) { Write-Output "FAILED: plied to Environmental StudiesApplication of geosynthetics i" }
$total += 1
$r = $d.Content.Find.Execute("tics in environmental projectsExperimental classes: soil mec", $true, $false, $false, $false, $false, $true, 1, $false, "tics in environmental projects^lExperimental classes: soil mec", 2)
if (-not $r) { Write-Output "FAILED: tics in environmental projectsExperimental classes: soil mec" }
$total += 1
$r = $d.Content.Find.Execute("ses: soil mechanics parametersThe discipline may have didact", $true, $false, $false, $false, $false, $true, 1, $false, "ses: soil mechanics parameters^lThe discipline may have didact", 2)
if (-not $r) { Write-Output "FAILED: ses: soil mechanics parametersThe discipline may have didact" }
$total += 1

# --- Paragraph 3 (Bibliografia) ---
$r = $d.Content.Find.Execute("Bibliografia básica:CHIOSSI, n. Geologia de Engenh", $true, $false, $false, $false, $false, $true, 1, $false, "Bibliografia básica:^lCHIOSSI, n. Geologia de Engenh", 2)
if (-not $r) { Write-Output "FAILED: Bibliografia básica:CHIOSSI, n. Geologia de Engenh" }
$total += 1
$r = $d.Content.Find.Execute("Textos: São Paulo, 1979. 427p.Bibliografia complementar", $true, $false, $false, $false, $false, $true, 1, $false, "Textos: São Paulo, 1979. 427p.^lBibliografia complementar", 2)
if (-not $r) { Write-Output "FAILED: Textos: São Paulo, 1979. 427p.Bibliografia complementar" }
$total += 1
$r = $d.Content.Find.Execute("Bibliografia complementarBOSCOV, M. E. Geotecnia ambien", $true, $false, $false, $false, $false, $true, 1, $false, "Bibliografia complementar^lBOSCOV, M. E. Geotecnia ambien", 2)
if (-not $r) { Write-Output "FAILED: Bibliografia complementarBOSCOV, M. E. Geotecnia ambien" }
$total += 1
$r = $d.Content.Find.Execute("Textos: São Paulo, 2008. 248p.ZUQUETTE, L. V. Geotecnia ambi", $true, $false, $false, $false, $false, $true, 1, $false, "Textos: São Paulo, 2008. 248p.^lZUQUETTE, L. V. Geotecnia ambi", 2)
if (-not $r) { Write-Output "FAILED: Textos: São Paulo, 2008. 248p.ZUQUETTE, L. V. Geotecnia ambi" }
$total += 1
$r = $d.Content.Find.Execute("r: Rio de Janeiro, 2015. 432p.OLIVEIRA, A. M. S.; JERÔNIMO, ", $true, $false, $false, $false, $false, $true, 1, $false, "r: Rio de Janeiro, 2015. 432p.^lOLIVEIRA, A. M. S.; JERÔNIMO, ", 2)
if (-not $r) { Write-Output "FAILED: r: Rio de Janeiro, 2015. 432p.OLIVEIRA, A. M. S.; JERÔNIMO, " }
$total += 1
$r = $d.Content.Find.Execute(" Ambiental, ABGE, 2018. 912 p.BARROW, C. J. Land degradation", $true, $false, $false, $false, $false, $true, 1, $false, " Ambiental, ABGE, 2018. 912 p.^lBARROW, C. J. Land degradation", 2)
if (-not $r) { Write-Output "FAILED:  Ambiental, ABGE, 2018. 912 p.BARROW, C. J. Land degradation" }
$total += 1
$r = $d.Content.Find.Execute("y Press: New York, 1991. 316p.KHALED, S. B. D. Fundamentos d", $true, $false, $false, $false, $false, $true, 1, $false, "y Press: New York, 1991. 316p.^lKHALED, S. B. D. Fundamentos d", 2)
if (-not $r) { Write-Output "FAILED: y Press: New York, 1991. 316p.KHALED, S. B. D. Fundamentos d" }
$total += 1
$r = $d.Content.Find.Execute(" Cengage Learning. 2020, 630p.SHARMA, H. D.; REDDY, K. R. Ge", $true, $false, $false, $false, $false, $true, 1, $false, " Cengage Learning. 2020, 630p.^lSHARMA, H. D.; REDDY, K. R. Ge", 2)
if (-not $r) { Write-Output "FAILED:  Cengage Learning. 2020, 630p.SHARMA, H. D.; REDDY, K. R. Ge" }
$total += 1
$r = $d.Content.Find.Execute("gineering, Wiley, 2004. 992p. YONG, R. N. Geoenvironmental e", $true, $false, $false, $false, $false, $true, 1, $false, "gineering, Wiley, 2004. 992p. ^lYONG, R. N. Geoenvironmental e", 2)
if (-not $r) { Write-Output "FAILED: gineering, Wiley, 2004. 992p. YONG, R. N. Geoenvironmental e" }
$total += 1
$r = $d.Content.Find.Execute("ration. CRC Press, 2001. 307p.DANIEL, D.E.  Geotechnical pra", $true, $false, $false, $false, $false, $true, 1, $false, "ration. CRC Press, 2001. 307p.^lDANIEL, D.E.  Geotechnical pra", 2)
if (-not $r) { Write-Output "FAILED: ration. CRC Press, 2001. 307p.DANIEL, D.E.  Geotechnical pra" }
$total += 1
$r = $d.Content.Find.Execute("l. Chapman & Hall, 1993. 693p.ROWE, R. K.; QUIGLEY, R.M.; BO", $true, $false, $false, $false, $false, $true, 1, $false, "l. Chapman & Hall, 1993. 693p.^lROWE, R. K.; QUIGLEY, R.M.; BO", 2)
if (-not $r) { Write-Output "FAILED: l. Chapman & Hall, 1993. 693p.ROWE, R. K.; QUIGLEY, R.M.; BO" }
$total += 1

Write-Output "Total replacements attempted: $total"